# The workbook contains two sheets that both had a header cell (B1) with a
# mislabeled value "ceRNAs". This upload corrects the typo to "ncRNAs" on
# both "Modules using ENCORI" and "Modules using SPONGEdb" sheets.
#
# Because both sheets end up sharing the exact same corrected string, the
# shared-strings table naturally collapses the two separate "ceRNAs" entries
# into a single "ncRNAs" entry (unique string count drops from 22 to 21),
# which is exactly what is reflected in the workbook's xl/sharedStrings.xml.

$wb = $excel.ActiveWorkbook

$wsEncori = $wb.Worksheets.Item("Modules using ENCORI")
$wsSponge = $wb.Worksheets.Item("Modules using SPONGEdb")

$wsEncori.Range("B1").Value = "ncRNAs"
$wsSponge.Range("B1").Value = "ncRNAs"
